$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column B values
$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 5

# Update A3 value (4 -> 2)
$ws.Range("A3").Value = 2

# Move the selection to B5 (matches the new <selection activeCell="B5" sqref="B5"/>)
$ws.Range("B5").Select()
